# "Wireframes version 2" -> "Wireframes version 1." revert.
# The paragraph currently reads "Version 2." built from runs:
#   "Versi" | "on" | " 2" | <bookmark _GoBack> | "."
# The target reads "Version 1." built from runs:
#   "Version" | " 1." | <bookmark _GoBack>
# We edit the live Range text (not the raw runs) and let Word's own
# run-(re)write logic collapse/split runs the same way a human typist
# retyping "Version" and then "1." would.

$d = $word.ActiveDocument

# --- 1) Merge "Versi" + "on" into a single run reading "Version" ---
$word1 = $d.Content.Duplicate
$word1.Find.ClearFormatting()
$word1.Find.Execute("Version", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Force an actual text change so Word rewrites the split runs as one run:
# first replace with a temporary different string, then retype the real word.
$word1.Text = "Versionx"
$word1b = $d.Range($word1.Start, $word1.Start + 8)
$word1b.Text = "Version"

# --- 2) Turn " 2" into " 1." ---
$num = $d.Content.Duplicate
$num.Find.ClearFormatting()
$num.Find.Execute(" 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$num.Text = " 1."

# --- 3) Remove the now-redundant trailing "." run (right after the bookmark) ---
$tail = $d.Content.Duplicate
$tail.Find.ClearFormatting()
$tail.Find.Execute("1..", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
# $tail now spans " 1.." ... narrow it to just the last character.
$dot = $d.Range($tail.End - 1, $tail.End)
$dot.Text = ""
